$d = $word.ActiveDocument

# The document currently has no headers. Add a default (primary) header to
# the first (only) section containing the questionnaire number, centered,
# in Arial 12pt - matching how Word materializes a brand-new header when a
# user types into it for the first time.
$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)

# Using Range.InsertAfter (rather than Range.Text = ...) only materializes
# the single "default" header story that we touch, instead of all six
# possible header/footer slots (default/even/first for both header and
# footer), which keeps the resulting package minimal.
$header.Range.InsertAfter("Questionnaire 30")

# Apply the paragraph-level formatting (style + centering) before doing any
# run-level formatting, otherwise setting Style afterwards clobbers the
# direct character formatting already applied to the run.
$header.Range.Paragraphs.Item(1).Style = "Header"
$header.Range.Paragraphs.Item(1).Alignment = 1

# Format only the text itself (excluding the trailing paragraph mark) so the
# run gets direct Arial/24-half-points (12pt) formatting without also
# stamping paragraph-mark run properties (w:pPr/w:rPr) that shouldn't be
# present.
$textRange = $header.Range.Duplicate()
$textRange.MoveEnd(1, -1) | Out-Null
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
